$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 16833.445
$ws.Range("I76").Value = 33975.25
$ws.Range("K76").Value = 33975.25
$ws.Range("M76").Value = -33660.25
$ws.Range("H79").Value = 16833.445
$ws.Range("I79").Value = 33975.25
$ws.Range("K79").Value = 33975.25
$ws.Range("M79").Value = -32883.25
$ws.Range("H128").Value = 33500
$ws.Range("J128").Value = 33500
$ws.Range("L128").Value = 33500
$ws.Range("N128").Value = -43460

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1424.625
$ws.Range("I61").Value = 1490.0322
$ws.Range("J61").Value = 1199.3334
$ws.Range("K61").Value = 1490.0322
$ws.Range("L61").Value = 1199.3334
$ws.Range("M61").Value = -1278.0322
$ws.Range("N61").Value = -1623.3334
$ws.Range("H136").Value = 1424.625
$ws.Range("I136").Value = 1490.0322
$ws.Range("J136").Value = 1199.3334
$ws.Range("K136").Value = 4470.096600000001
$ws.Range("L136").Value = 3598.0002
$ws.Range("M136").Value = -1920.096600000001
$ws.Range("N136").Value = -8698.0002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 85236.836
$ws.Range("I134").Value = 108932.18
$ws.Range("J134").Value = 2303.125
$ws.Range("K134").Value = 326796.54
$ws.Range("L134").Value = 6909.375
$ws.Range("M134").Value = -324261.54
$ws.Range("N134").Value = -11979.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2359.1428
$ws.Range("I99").Value = 1975
$ws.Range("J99").Value = 2871.3333
$ws.Range("K99").Value = 1975
$ws.Range("L99").Value = 2871.3333
$ws.Range("M99").Value = -477
$ws.Range("N99").Value = -5867.3333
$ws.Range("H100").Value = 54000
$ws.Range("J100").Value = 54000
$ws.Range("L100").Value = 54000
$ws.Range("N100").Value = -56164
$ws.Range("H126").Value = 2359.1428
$ws.Range("I126").Value = 1975
$ws.Range("J126").Value = 2871.3333
$ws.Range("K126").Value = 5925
$ws.Range("L126").Value = 8613.999899999999
$ws.Range("M126").Value = -3455
$ws.Range("N126").Value = -13553.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 278137.97
$ws.Range("I5").Value = 298.0263
$ws.Range("J5").Value = 1333929.8
$ws.Range("K5").Value = 894.0789
$ws.Range("L5").Value = 4001789.4
$ws.Range("M5").Value = -782.0789
$ws.Range("N5").Value = -4002013.4
$ws.Range("H124").Value = 5485.7144
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 5485.7144
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 16457.1432
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -26277.1432
$ws.Range("H125").Value = 2583.2258
$ws.Range("I125").Value = 1160
$ws.Range("J125").Value = 2735.7144
$ws.Range("K125").Value = 3480
$ws.Range("L125").Value = 8207.143199999999
$ws.Range("M125").Value = 1440
$ws.Range("N125").Value = -18047.1432
$ws.Range("H126").Value = 1922.5807
$ws.Range("J126").Value = 1922.5807
$ws.Range("L126").Value = 5767.742099999999
$ws.Range("N126").Value = -15647.7421
$ws.Range("H131").Value = 953.3333
$ws.Range("I131").Value = 504.44446
$ws.Range("J131").Value = 1017.4603
$ws.Range("K131").Value = 1513.33338
$ws.Range("L131").Value = 3052.3809
$ws.Range("M131").Value = 3526.66662
$ws.Range("N131").Value = -13132.3809
$ws.Range("H135").Value = 278137.97
$ws.Range("I135").Value = 298.0263
$ws.Range("J135").Value = 1333929.8
$ws.Range("K135").Value = 2682.2367
$ws.Range("L135").Value = 12005368.2
$ws.Range("M135").Value = -147.2366999999999
$ws.Range("N135").Value = -12010438.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1776.6666
$ws.Range("I122").Value = 1686.7142
$ws.Range("J122").Value = 1956.5714
$ws.Range("K122").Value = 5060.142599999999
$ws.Range("L122").Value = 5869.7142
$ws.Range("M122").Value = -2610.142599999999
$ws.Range("N122").Value = -10769.7142
$ws.Range("H125").Value = 29666.666
$ws.Range("J125").Value = 29666.666
$ws.Range("L125").Value = 29666.666
$ws.Range("N125").Value = -34586.666
$ws.Range("H126").Value = 4978.9653
$ws.Range("I126").Value = 2892.1428
$ws.Range("J126").Value = 6926.6665
$ws.Range("K126").Value = 8676.428400000001
$ws.Range("L126").Value = 20779.9995
$ws.Range("M126").Value = -6206.428400000001
$ws.Range("N126").Value = -25719.9995
$ws.Range("H134").Value = 13238
$ws.Range("J134").Value = 13238
$ws.Range("L134").Value = 39714
$ws.Range("N134").Value = -44784

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7145245
$ws.Range("I7").Value = 12501914
$ws.Range("J7").Value = 3019.8333
$ws.Range("K7").Value = 12501914
$ws.Range("L7").Value = 3019.8333
$ws.Range("M7").Value = -12501802
$ws.Range("N7").Value = -3243.8333
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H40").Value = 1247.28
$ws.Range("I40").Value = 1259
$ws.Range("J40").Value = 1161.3334
$ws.Range("K40").Value = 1259
$ws.Range("L40").Value = 1161.3334
$ws.Range("M40").Value = -1123
$ws.Range("N40").Value = -1433.3334
$ws.Range("H126").Value = 7145245
$ws.Range("I126").Value = 12501914
$ws.Range("J126").Value = 3019.8333
$ws.Range("K126").Value = 37505742
$ws.Range("L126").Value = 9059.499899999999
$ws.Range("M126").Value = -37503272
$ws.Range("N126").Value = -13999.4999
$ws.Range("H135").Value = 26249.75
$ws.Range("J135").Value = 26249.75
$ws.Range("L135").Value = 26249.75
$ws.Range("N135").Value = -36389.75
$ws.Range("H136").Value = 1226.2742
$ws.Range("I136").Value = 1076.4814
$ws.Range("J136").Value = 2237.375
$ws.Range("K136").Value = 3229.4442
$ws.Range("L136").Value = 6712.125
$ws.Range("M136").Value = -679.4441999999999
$ws.Range("N136").Value = -11812.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2173.5
$ws.Range("I126").Value = 1745.4783
$ws.Range("J126").Value = 2930.7693
$ws.Range("K126").Value = 5236.4349
$ws.Range("L126").Value = 8792.3079
$ws.Range("M126").Value = -2766.4349
$ws.Range("N126").Value = -13732.3079
$ws.Range("I132").Value = 2347.973
$ws.Range("J132").Value = 824.6667
$ws.Range("K132").Value = 7043.919
$ws.Range("L132").Value = 2474.0001
$ws.Range("M132").Value = -4513.919
$ws.Range("N132").Value = -7534.0001
